$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: rename placeholder '5' -> '6' and refresh all financial metrics; drop stale debt_ebitda/net_debt_ebitda
$ws.Range("B2").Value = "'6"
$ws.Range("B2").Style = "Normal"  # keep plain formatting; value is text "6", not a number
$ws.Range("D2").Value = 0.07735
$ws.Range("E2").Value = -0.0201
$ws.Range("F2").Value = 0.0627
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1058.2
$ws.Range("L2").Value = 0.1284831412926021
$ws.Range("M2").Value = 1185.001
$ws.Range("N2").Value = 0.03534035566981101
$ws.Range("O2").Value = 1.119827064827065
$ws.Range("P2").Value = 1185
$ws.Range("Q2").Value = 0.03534032584675122
$ws.Range("R2").Value = 1.11982611982612
$ws.Range("S2").Value = 0.0009999999999976694
$ws.Range("T2").Value = ([double]"8.438811444021309e-07")
$ws.Range("U2").Value = 14258.2
$ws.Range("V2").Value = 0.4252231510448509
$ws.Range("W2").Value = 0.09603758821697028
$ws.Range("X2").Value = 0.1165098762020992
$ws.Range("Y2").Value = -0.02047228798512891
$ws.Range("Z2").Value = 0.07171385856102092
$ws.Range("AB2").Value = 0.03806368920158379
$ws.Range("AC2").Value = -0.03806368920158379
$ws.Range("AD2").Value = 128845.5
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 128845.5
$ws.Range("AG2").Value = 114587.3
$ws.Range("AH2").Value = 0.7934979547545643
$ws.Range("AI2").Value = 0.8564531696188801
$ws.Range("AJ2").Value = 0.7736196178192581
$ws.Range("AK2").Value = 0.8414239410908588
$ws.Range("AN2").ClearContents()
$ws.Range("AP2").ClearContents()

# Row 3: Grupo Security S.A. -- refresh financial metrics
$ws.Range("D3").Value = 0.0713
$ws.Range("E3").Value = 0.0169
$ws.Range("K3").Value = 86.3
$ws.Range("L3").Value = 0.1870394451668834
$ws.Range("M3").Value = 57.601
$ws.Range("N3").Value = 0.0772442000804613
$ws.Range("O3").Value = 0.6674507531865586
$ws.Range("P3").Value = 57.6
$ws.Range("Q3").Value = 0.07724285905860265
$ws.Range("R3").Value = 0.6674391657010429
$ws.Range("S3").Value = 0.0009999999999976694
$ws.Range("T3").Value = ([double]"1.736080970812433e-05")
$ws.Range("U3").Value = 743.6
$ws.Range("V3").Value = 0.9971838540968218
$ws.Range("W3").Value = 0.08446706469609473
$ws.Range("X3").Value = 0.197259728893272
$ws.Range("Y3").Value = -0.1127926641971773
$ws.Range("Z3").Value = 0.08382841881506514
$ws.Range("AB3").Value = 0.03721833498934456
$ws.Range("AC3").Value = -0.03721833498934456
$ws.Range("AD3").Value = 6394.7
$ws.Range("AF3").Value = 6394.7
$ws.Range("AG3").Value = 5651.099999999999
$ws.Range("AH3").Value = 0.8955660747297071
$ws.Range("AI3").Value = 0.8537080301715506
$ws.Range("AJ3").Value = 0.883426088044022
$ws.Range("AK3").Value = 0.8375846685144288

# Row 4: Banco de Chile -- refresh financial metrics; drop stale debt_ebitda/net_debt_ebitda
$ws.Range("D4").Value = 0.0196
$ws.Range("E4").Value = -0.0242
$ws.Range("F4").Value = 0.061
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 615.1
$ws.Range("L4").Value = 0.3204647285610087
$ws.Range("M4").Value = 445.5
$ws.Range("N4").Value = 0.04320377051088095
$ws.Range("O4").Value = 0.7242724760201593
$ws.Range("P4").Value = 445.5
$ws.Range("Q4").Value = 0.04320377051088095
$ws.Range("R4").Value = 0.7242724760201593
$ws.Range("U4").Value = 2001.7
$ws.Range("V4").Value = 0.1941211839093836
$ws.Range("W4").Value = 0.1308278033010039
$ws.Range("X4").Value = 0.07271681589603375
$ws.Range("Y4").Value = 0.05811098740497016
$ws.Range("Z4").Value = 0.08904910366328918
$ws.Range("AA4").Value = 0
$ws.Range("AB4").Value = 0.03761613689500742
$ws.Range("AC4").Value = -0.03761613689500742
$ws.Range("AD4").Value = 20475.1
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 20475.1
$ws.Range("AG4").Value = 18473.4
$ws.Range("AH4").Value = 0.6650631603906882
$ws.Range("AI4").Value = 0.8160466469514479
$ws.Range("AJ4").Value = 0.6417717561229807
$ws.Range("AK4").Value = 0.800098748749399
$ws.Range("AN4").ClearContents()
$ws.Range("AP4").ClearContents()

# Row 5: company re-mapped to Banco Santander-Chile; refresh financial metrics
$ws.Range("B5").Value = "Banco Santander-Chile (SNSE:BSANTANDER)"
$ws.Range("D5").Value = 0.0199
$ws.Range("E5").Value = -0.022
$ws.Range("F5").Value = 0.0644
$ws.Range("K5").Value = 572.9
$ws.Range("L5").Value = 0.3100443770970884
$ws.Range("M5").Value = 210.5
$ws.Range("N5").Value = 0.02328179264273232
$ws.Range("O5").Value = 0.3674288706580555
$ws.Range("P5").Value = 210.5
$ws.Range("Q5").Value = 0.02328179264273232
$ws.Range("R5").Value = 0.3674288706580555
$ws.Range("U5").Value = 3389
$ws.Range("V5").Value = 0.3748313314309731
$ws.Range("W5").Value = 0.1242302020990545
$ws.Range("X5").Value = 0.1027010192415178
$ws.Range("Y5").Value = 0.02152918285753673
$ws.Range("Z5").Value = 0.06836790948411421
$ws.Range("AB5").Value = 0.0380401924282074
$ws.Range("AC5").Value = -0.0380401924282074
$ws.Range("AD5").Value = 32297.3
$ws.Range("AF5").Value = 32297.3
$ws.Range("AG5").Value = 28908.3
$ws.Range("AH5").Value = 0.7812848493058563
$ws.Range("AI5").Value = 0.872030326512208
$ws.Range("AJ5").Value = 0.7617530573364216
$ws.Range("AK5").Value = 0.8591412837056696

# Row 6: company re-mapped to Scotiabank Chile, S.A.; refresh financial metrics
$ws.Range("B6").Value = "Scotiabank Chile, S.A. (SNSE:SCOTIABKCL)"
$ws.Range("D6").Value = 0.241
$ws.Range("E6").Value = 0.248
$ws.Range("F6").ClearContents()
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 320.5
$ws.Range("L6").Value = 0.2659751037344398
$ws.Range("M6").Value = 129.3
$ws.Range("N6").Value = 0.02202162990717875
$ws.Range("O6").Value = 0.4034321372854914
$ws.Range("P6").Value = 129.3
$ws.Range("Q6").Value = 0.02202162990717875
$ws.Range("R6").Value = 0.4034321372854914
$ws.Range("U6").Value = 2172.4
$ws.Range("V6").Value = 0.3699906327173635
$ws.Range("W6").Value = 0.1076081117378458
$ws.Range("X6").Value = 0.1081079468492021
$ws.Range("Y6").Value = -0.0004998351113562893
$ws.Range("Z6").Value = 0.06581318885381279
$ws.Range("AA6").Value = 0
$ws.Range("AB6").Value = 0.03808718597496018
$ws.Range("AC6").Value = -0.03808718597496018
$ws.Range("AD6").Value = 22653.7
$ws.Range("AE6").Value = 0
$ws.Range("AF6").Value = 22653.7
$ws.Range("AG6").Value = 20481.3
$ws.Range("AH6").Value = 0.7941644580931948
$ws.Range("AI6").Value = 0.8883629732750339
$ws.Range("AJ6").Value = 0.7771963510518807
$ws.Range("AK6").Value = 0.877966915436748
$ws.Range("AN6").ClearContents()
$ws.Range("AP6").ClearContents()

# Row 7: company re-mapped to Banco de Crédito e Inversiones (BCI); refresh financial metrics
$ws.Range("B7").Value = "Banco de Crédito e Inversiones (SNSE:BCI)"
$ws.Range("D7").Value = 0.0834
$ws.Range("E7").Value = -0.0201
$ws.Range("K7").Value = 387
$ws.Range("L7").Value = 0.2157792026763312
$ws.Range("M7").Value = 180
$ws.Range("N7").Value = 0.03083353318030765
$ws.Range("O7").Value = 0.4651162790697674
$ws.Range("P7").Value = 180
$ws.Range("Q7").Value = 0.03083353318030765
$ws.Range("R7").Value = 0.4651162790697674
$ws.Range("T7").Value = 0
$ws.Range("U7").Value = 3262.7
$ws.Range("V7").Value = 0.5588920483743876
$ws.Range("W7").Value = 0.07693378128541041
$ws.Range("X7").Value = 0.1249118055549963
$ws.Range("Y7").Value = -0.04797802426958588
$ws.Range("Z7").Value = 0.07286533219034774
$ws.Range("AB7").Value = 0.03820337058934214
$ws.Range("AC7").Value = -0.03820337058934214
$ws.Range("AD7").Value = 27714.2
$ws.Range("AF7").Value = 27714.2
$ws.Range("AG7").Value = 24451.5
$ws.Range("AH7").Value = 0.8260073915116833
$ws.Range("AI7").Value = 0.8464780730957894
$ws.Range("AJ7").Value = 0.8072652718946955
$ws.Range("AK7").Value = 0.8294858181892196

# Row 8 (new): Itaú Corpbanca -- add full record
$ws.Range("A8").Value = "Chile"
$ws.Range("B8").Value = "Itaú Corpbanca (SNSE:ITAUCORP)"
$ws.Range("C8").Value = "Bank (Money Center)"
$ws.Range("D8").Value = 0.239
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = -923.6
$ws.Range("L8").Value = -0.9153617443012885
$ws.Range("M8").Value = 162.1
$ws.Range("N8").Value = 0.09407463292902328
$ws.Range("O8").Value = -0.1755088783022954
$ws.Range("P8").Value = 162.1
$ws.Range("Q8").Value = 0.09407463292902328
$ws.Range("R8").Value = -0.1755088783022954
$ws.Range("S8").Value = 0
$ws.Range("T8").Value = 0
$ws.Range("U8").Value = 2688.8
$ws.Range("V8").Value = 1.560443386918925
$ws.Range("W8").Value = -0.1962600934976626
$ws.Range("X8").Value = 0.2469915834728775
$ws.Range("Y8").Value = -0.4432516769705401
$ws.Range("Z8").Value = 0.05656590572722788
$ws.Range("AA8").Value = 0
$ws.Range("AB8").Value = 0.04068734165018113
$ws.Range("AC8").Value = -0.04068734165018113
$ws.Range("AD8").Value = 19310.5
$ws.Range("AE8").Value = 0
$ws.Range("AF8").Value = 19310.5
$ws.Range("AG8").Value = 16621.7
$ws.Range("AH8").Value = 0.9180786931385974
$ws.Range("AI8").Value = 0.8551393384909018
$ws.Range("AJ8").Value = 0.9060714752954516
$ws.Range("AK8").Value = 0.8355594206978368
$ws.Range("AL8").Value = 0
$ws.Range("AM8").Value = 0
